# Updated calibration data with new costs
# Target worksheet is the only sheet in the workbook: "strategy_id-0"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (gdp_mmm_usd): replace the J2:AS2 series with newly calibrated values ---
$row2Values = @(
    1857.648102,
    1869.8082082,
    1881.9683144,
    1894.1284206,
    1906.2885268,
    1918.448633,
    1943.3463408,
    1968.2440486,
    1993.1417564,
    2018.0394642,
    2042.937172,
    2099.7419388,
    2156.5467056,
    2213.3514724,
    2270.1562392,
    2326.961006,
    2388.2460956,
    2449.5311852,
    2510.8162748,
    2572.1013644,
    2633.386454,
    2698.6293606,
    2763.8722672,
    2829.1151738,
    2894.3580804,
    2959.600987,
    3030.6637982,
    3101.7266094,
    3172.7894206,
    3243.8522318,
    3314.915043,
    3393.2952838,
    3471.6755246,
    3550.0557654,
    3628.4360062,
    3706.816247
)

$startCol = 10   # column J
for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $ws.Cells.Item(2, $startCol + $i).Value = $row2Values[$i]
}

# --- Row 8 (elasticity_gnrl_rate_occupancy_to_gdppc): constant -0.05 -> -0.1 across J8:AS8 ---
for ($col = 10; $col -le 45; $col++) {
    $ws.Cells.Item(8, $col).Value = -0.1
}

# --- Row 9 (frac_gnrl_eating_red_meat): constant 1.05404749105734 -> 1 across J9:AS9 ---
for ($col = 10; $col -le 45; $col++) {
    $ws.Cells.Item(9, $col).Value = 1
}

# --- Row 13 (occrateinit_gnrl_occupancy): constant 4.296989118339568 -> 3.145207224 across J13:AS13 ---
for ($col = 10; $col -le 45; $col++) {
    $ws.Cells.Item(13, $col).Value = 3.145207224
}
